$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.2750083333333334
$ws.Range("H2").Value = 0.825025
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.090355666666666
$ws.Range("N2").Value = 9.271066999999999
$ws.Range("O2").Value = 0.06928583878088775
$ws.Range("P2").Value = 0.06928583878088775
$ws.Range("Q2").Value = 0.8498735612972221
$ws.Range("R2").Value = 7.648862051674999
$ws.Range("S2").Value = 0.06928583878088775
$ws.Range("T2").Value = 0.06928583878088775

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.2750083333333334
$ws.Range("H3").Value = 0.825025
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 25.17096033333333
$ws.Range("N3").Value = 75.51288099999999
$ws.Range("O3").Value = 0.5643334579338453
$ws.Range("P3").Value = 0.5643334579338454
$ws.Range("Q3").Value = 6.922223849669444
$ws.Range("R3").Value = 62.30001464702499
$ws.Range("S3").Value = 0.5643334579338453
$ws.Range("T3").Value = 0.5643334579338454

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.2750083333333334
$ws.Range("H4").Value = 0.825025
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 16.34167533333333
$ws.Range("N4").Value = 49.025026
$ws.Range("O4").Value = 0.366380703285267
$ws.Range("P4").Value = 0.366380703285267
$ws.Range("Q4").Value = 4.494096897294444
$ws.Range("R4").Value = 40.44687207565
$ws.Range("S4").Value = 0.366380703285267
$ws.Range("T4").Value = 0.366380703285267
